# Update of 2025 data and RF changes
# Change the RF (column I) values for rows 22 through 51 from 2.357424242424242
# to 4.17304347826087 on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I22:I51").Value = 4.17304347826087
